$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add "LE PAYS" sheet right after GÉNÉRALITÉS
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = 'LE PAYS'

# Add "LE RELIEF ET LA VÉGÉTATION" sheet right after LE PAYS
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = 'LE RELIEF ET LA VÉGÉTATION'

# Populate sheet2
$ws2.Range("A1").Value = 'Mot français'
$ws2.Range("B1").Value = 'Grammaire'
$ws2.Range("C1").Value = 'Prononciation'
$ws2.Range("D1").Value = 'Signification en tchèque'
$ws2.Range("A2").Value = 'capitale'
$ws2.Range("B2").Value = 'nf'
$ws2.Range("C2").Value = 'kapital'
$ws2.Range("D2").Value = 'hlavní město'
$ws2.Range("F2").Formula = '= "{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"'
$ws2.Range("A3").Value = 'frontière'
$ws2.Range("B3").Value = 'nf'
$ws2.Range("C3").Value = 'fro~tje:r'
$ws2.Range("D3").Value = 'hranice'
$ws2.Range("F3").Formula = '= "{ ""foreign"": """ & A3 & """, ""grammar"": """ & B3 & """, ""pronunciation"": """ & C3 & """, ""meaning"": """ & D3 & """ },"'
$ws2.Range("A4").Value = 'indigène'
$ws2.Range("B4").Value = 'adj'
$ws2.Range("C4").Value = 'e~dižen'
$ws2.Range("D4").Value = 'domorodý, tuzemský, domácí'
$ws2.Range("F4").Formula = '= "{ ""foreign"": """ & A4 & """, ""grammar"": """ & B4 & """, ""pronunciation"": """ & C4 & """, ""meaning"": """ & D4 & """ },"'
$ws2.Range("A5").Value = 'peuplement'
$ws2.Range("B5").Value = 'nm'
$ws2.Range("C5").Value = 'pöpl@ma~'
$ws2.Range("D5").Value = 'zalidnění'
$ws2.Range("F5").Formula = '= "{ ""foreign"": """ & A5 & """, ""grammar"": """ & B5 & """, ""pronunciation"": """ & C5 & """, ""meaning"": """ & D5 & """ },"'
$ws2.Range("A6").Value = 'peupler'
$ws2.Range("B6").Value = 'vt'
$ws2.Range("C6").Value = 'pöple.'
$ws2.Range("D6").Value = 'zalidnit, obývat'
$ws2.Range("F6").Formula = '= "{ ""foreign"": """ & A6 & """, ""grammar"": """ & B6 & """, ""pronunciation"": """ & C6 & """, ""meaning"": """ & D6 & """ },"'
$ws2.Range("A7").Value = 'région'
$ws2.Range("B7").Value = 'nf'
$ws2.Range("C7").Value = 're.žjo~'
$ws2.Range("D7").Value = 'oblast'
$ws2.Range("F7").Formula = '= "{ ""foreign"": """ & A7 & """, ""grammar"": """ & B7 & """, ""pronunciation"": """ & C7 & """, ""meaning"": """ & D7 & """ },"'
$ws2.Range("A8").Value = 'territoire'
$ws2.Range("B8").Value = 'nm'
$ws2.Range("C8").Value = 'teritu^a:r'
$ws2.Range("D8").Value = 'území'
$ws2.Range("F8").Formula = '= "{ ""foreign"": """ & A8 & """, ""grammar"": """ & B8 & """, ""pronunciation"": """ & C8 & """, ""meaning"": """ & D8 & """ },"'

# Populate sheet3
$ws3.Range("A1").Value = 'Mot français'
$ws3.Range("B1").Value = 'Grammaire'
$ws3.Range("C1").Value = 'Prononciation'
$ws3.Range("D1").Value = 'Signification en tchèque'
$ws3.Range("A2").Value = 'désert'
$ws3.Range("B2").Value = 'nm'
$ws3.Range("C2").Value = 'de.ze:r'
$ws3.Range("D2").Value = 'pustina, poušť'
$ws3.Range("F2").Formula = '= "{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"'
$ws3.Range("A3").Value = 'désertique'
$ws3.Range("B3").Value = 'adj'
$ws3.Range("C3").Value = 'de.zertik'
$ws3.Range("D3").Value = 'pouštní'
$ws3.Range("F3").Formula = '= "{ ""foreign"": """ & A3 & """, ""grammar"": """ & B3 & """, ""pronunciation"": """ & C3 & """, ""meaning"": """ & D3 & """ },"'
$ws3.Range("A4").Value = 'forestier, ière'
$ws3.Range("B4").Value = 'adj'
$ws3.Range("C4").Value = 'forestje., stje:r'
$ws3.Range("D4").Value = 'lesní, lesnický'
$ws3.Range("F4").Formula = '= "{ ""foreign"": """ & A4 & """, ""grammar"": """ & B4 & """, ""pronunciation"": """ & C4 & """, ""meaning"": """ & D4 & """ },"'
$ws3.Range("A5").Value = 'forêt'
$ws3.Range("B5").Value = 'nf'
$ws3.Range("C5").Value = 'fore'
$ws3.Range("D5").Value = 'les'
$ws3.Range("F5").Formula = '= "{ ""foreign"": """ & A5 & """, ""grammar"": """ & B5 & """, ""pronunciation"": """ & C5 & """, ""meaning"": """ & D5 & """ },"'
$ws3.Range("A6").Value = 'jungle'
$ws3.Range("B6").Value = 'nf'
$ws3.Range("C6").Value = 'žo~:gl, žö~:gl'
$ws3.Range("D6").Value = 'džungle'
$ws3.Range("F6").Formula = '= "{ ""foreign"": """ & A6 & """, ""grammar"": """ & B6 & """, ""pronunciation"": """ & C6 & """, ""meaning"": """ & D6 & """ },"'
$ws3.Range("A7").Value = 'oasis'
$ws3.Range("B7").Value = 'nf'
$ws3.Range("C7").Value = 'oazis'
$ws3.Range("D7").Value = 'oáza'
$ws3.Range("F7").Formula = '= "{ ""foreign"": """ & A7 & """, ""grammar"": """ & B7 & """, ""pronunciation"": """ & C7 & """, ""meaning"": """ & D7 & """ },"'
$ws3.Range("A8").Value = 'plaine'
$ws3.Range("B8").Value = 'nf'
$ws3.Range("C8").Value = 'plen'
$ws3.Range("D8").Value = 'rovina, planina'
$ws3.Range("F8").Formula = '= "{ ""foreign"": """ & A8 & """, ""grammar"": """ & B8 & """, ""pronunciation"": """ & C8 & """, ""meaning"": """ & D8 & """ },"'
$ws3.Range("A9").Value = 'plateau'
$ws3.Range("B9").Value = 'nm'
$ws3.Range("C9").Value = 'plato.'
$ws3.Range("D9").Value = 'plošina, náhorní rovina'
$ws3.Range("F9").Formula = '= "{ ""foreign"": """ & A9 & """, ""grammar"": """ & B9 & """, ""pronunciation"": """ & C9 & """, ""meaning"": """ & D9 & """ },"'
$ws3.Range("A10").Value = 'pôle'
$ws3.Range("B10").Value = 'nm'
$ws3.Range("C10").Value = 'po:l'
$ws3.Range("D10").Value = 'pól'
$ws3.Range("F10").Formula = '= "{ ""foreign"": """ & A10 & """, ""grammar"": """ & B10 & """, ""pronunciation"": """ & C10 & """, ""meaning"": """ & D10 & """ },"'
$ws3.Range("A11").Value = 'prairie'
$ws3.Range("B11").Value = 'nf'
$ws3.Range("C11").Value = 'preri'
$ws3.Range("D11").Value = 'louka, prérie'
$ws3.Range("F11").Formula = '= "{ ""foreign"": """ & A11 & """, ""grammar"": """ & B11 & """, ""pronunciation"": """ & C11 & """, ""meaning"": """ & D11 & """ },"'
$ws3.Range("A12").Value = 'rural, ale, aux'
$ws3.Range("B12").Value = 'adj'
$ws3.Range("C12").Value = 'rüral'
$ws3.Range("D12").Value = 'polní, venkovský'
$ws3.Range("F12").Formula = '= "{ ""foreign"": """ & A12 & """, ""grammar"": """ & B12 & """, ""pronunciation"": """ & C12 & """, ""meaning"": """ & D12 & """ },"'
$ws3.Range("A13").Value = 'savane'
$ws3.Range("B13").Value = 'nf'
$ws3.Range("C13").Value = 'savan'
$ws3.Range("D13").Value = 'savana'
$ws3.Range("F13").Formula = '= "{ ""foreign"": """ & A13 & """, ""grammar"": """ & B13 & """, ""pronunciation"": """ & C13 & """, ""meaning"": """ & D13 & """ },"'
$ws3.Range("A14").Value = 'sol'
$ws3.Range("B14").Value = 'nm'
$ws3.Range("C14").Value = 'sol'
$ws3.Range("D14").Value = 'půda, země'
$ws3.Range("F14").Formula = '= "{ ""foreign"": """ & A14 & """, ""grammar"": """ & B14 & """, ""pronunciation"": """ & C14 & """, ""meaning"": """ & D14 & """ },"'
$ws3.Range("A15").Value = 'steppe'
$ws3.Range("B15").Value = 'nf'
$ws3.Range("C15").Value = 'step'
$ws3.Range("D15").Value = 'step'
$ws3.Range("F15").Formula = '= "{ ""foreign"": """ & A15 & """, ""grammar"": """ & B15 & """, ""pronunciation"": """ & C15 & """, ""meaning"": """ & D15 & """ },"'

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 37.85546875
$ws2.Columns.Item(2).ColumnWidth = 11
$ws2.Columns.Item(3).ColumnWidth = 34
$ws2.Columns.Item(4).ColumnWidth = 63

$ws3.Columns.Item(1).ColumnWidth = 36.28515625
$ws3.Columns.Item(2).ColumnWidth = 11
$ws3.Columns.Item(3).ColumnWidth = 27.42578125
$ws3.Columns.Item(4).ColumnWidth = 47.5703125

# Selections / active cell per sheet
$ws1.Activate()
$ws1.Range("F2").Select()
$ws2.Activate()
$ws2.Range("F2").Select()
$ws3.Activate()
$ws3.Range("F2").Select()

